$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.578.96'
$ws.Range("E2").Value = '  +0.71%  '

$ws.Range("D3").Value = '3.390.39'
$ws.Range("E3").Value = '  +0.47%  '

$ws.Range("E4").Value = '  +0.09%  '

$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '576.29'
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = '  +0.58%  '

$origStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '136.76'
$ws.Range("D6").Style = $origStyle
$ws.Range("E6").Value = '  +0.54%  '

$ws.Range("E7").Value = '  -0.02%  '

$ws.Range("D8").Value = '3.390.97'
$ws.Range("E8").Value = '  +0.53%  '

$ws.Range("E9").Value = '  -0.64%  '

$origStyle = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.50'
$ws.Range("D10").Style = $origStyle
$ws.Range("E10").Value = '  -1.44%  '

$origStyle = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.125'
$ws.Range("D11").Style = $origStyle
$ws.Range("E11").Value = '  +2.16%  '

$ws.Range("E12").Value = '  +0.07%  '

$ws.Range("D13").Value = '3.967.08'
$ws.Range("E13").Value = '  +0.55%  '

$ws.Range("E14").Value = '  +1.37%  '

$ws.Range("E15").Value = '  +1.69%  '

$ws.Range("D16").Value = '3.393.42'
$ws.Range("E16").Value = '  +0.60%  '

$origStyle = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '25.70'
$ws.Range("D17").Style = $origStyle
$ws.Range("E17").Value = '  +1.87%  '

$ws.Range("D18").Value = '61.688.97'
$ws.Range("E18").Value = '  +0.67%  '

$ws.Range("E19").Value = '  +1.56%  '

$ws.Range("B20").Value = 'Polkadot'
$ws.Range("C20").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$origStyle = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.86'
$ws.Range("D20").Style = $origStyle
$ws.Range("E20").Value = '  +0.76%  '

$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$origStyle = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.48'
$ws.Range("D21").Style = $origStyle
$ws.Range("E21").Value = '  +0.24%  '

$origStyle = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '379.37'
$ws.Range("D22").Style = $origStyle
$ws.Range("E22").Value = '  +1.15%  '

$origStyle = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.559'
$ws.Range("D23").Style = $origStyle
$ws.Range("E23").Value = '  -1.48%  '

$ws.Range("D24").Value = '3.525.39'
$ws.Range("E24").Value = '  +0.42%  '

$ws.Range("E25").Value = '  +0.02%  '

$origStyle = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000124'
$ws.Range("D26").Style = $origStyle
$ws.Range("E26").Value = '  +6.49%  '

$origStyle = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '71.12'
$ws.Range("D27").Style = $origStyle
$ws.Range("E27").Value = '  +0.90%  '

$origStyle = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.69'
$ws.Range("D28").Style = $origStyle
$ws.Range("E28").Value = '  +0.69%  '

$origStyle = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.57'
$ws.Range("D29").Style = $origStyle
$ws.Range("E29").Value = '  -1.79%  '

$ws.Range("E30").Value = '  +0.58%  '

$ws.Range("E31").Value = '  +3.23%  '

$origStyle = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.19'
$ws.Range("D32").Style = $origStyle
$ws.Range("E32").Value = '  +0.63%  '

$ws.Range("E33").Value = '  +0.75%  '

$origStyle = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '23.38'
$ws.Range("D35").Style = $origStyle
$ws.Range("E35").Value = '  +0.04%  '

$ws.Range("E36").Value = '  -4.09%  '

$ws.Range("E37").Value = '  -0.05%  '

$origStyle = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.82'
$ws.Range("D38").Style = $origStyle
$ws.Range("E38").Value = '  -1.47%  '

$origStyle = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '164.70'
$ws.Range("D39").Style = $origStyle
$ws.Range("E39").Value = '  +0.15%  '

$origStyle = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0781'
$ws.Range("D40").Style = $origStyle
$ws.Range("E40").Value = '  -1.05%  '

$ws.Range("E41").Value = '  +2.58%  '

$ws.Range("E42").Value = '  +2.72%  '

$ws.Range("E43").Value = '  +0.07%  '

$ws.Range("E44").Value = '  +7.44%  '

$origStyle = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '25.19'
$ws.Range("D45").Style = $origStyle
$ws.Range("E45").Value = '  +8.31%  '

$ws.Range("E46").Value = '  -0.26%  '

$origStyle = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '41.38'
$ws.Range("D47").Style = $origStyle
$ws.Range("E47").Value = '  +0.05%  '

$origStyle = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '6.85'
$ws.Range("D48").Style = $origStyle
$ws.Range("E48").Value = '  -1.75%  '

$origStyle = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '22.71'
$ws.Range("D49").Style = $origStyle
$ws.Range("E49").Value = '  -0.53%  '

$ws.Range("D50").Value = '2.344.20'
$ws.Range("E50").Value = '  +5.69%  '

$ws.Range("E51").Value = '  +1.88%  '
